$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep exact text formatting (leading/trailing zeros, etc.)
$priceCells = @("D2", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price values
$ws.Range("D2").Value = "246.57"
$ws.Range("D4").Value = "5.421"
$ws.Range("D5").Value = "0.05784"
$ws.Range("D6").Value = "3.385"
$ws.Range("D7").Value = "6.335"
$ws.Range("D8").Value = "0.8081"
$ws.Range("D9").Value = "0.9458"
$ws.Range("D11").Value = "0.07500"
$ws.Range("D12").Value = "0.03196"
$ws.Range("D13").Value = "0.03019"
$ws.Range("D14").Value = "4.167"
$ws.Range("D15").Value = "0.09403"
$ws.Range("D16").Value = "0.001594"
$ws.Range("D17").Value = "0.04814"
$ws.Range("D18").Value = "0.0005888"
$ws.Range("D19").Value = "0.006181"
$ws.Range("D20").Value = "0.004110"
$ws.Range("D21").Value = "0.0009983"
$ws.Range("D23").Value = "3.772"
$ws.Range("D24").Value = "2.232"
$ws.Range("D26").Value = "0.1260"
$ws.Range("D27").Value = "0.0002856"
$ws.Range("D40").Value = "0.03893"
$ws.Range("D41").Value = "0.006328"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").Value = "0.002999"
$ws.Range("D44").Value = "0.006338"
$ws.Range("D45").Value = "0.00005587"
$ws.Range("D48").Value = "0.1436"

# Other text / numeric updates
$ws.Range("G2").Value = "7"
$ws.Range("G3").Value = "7"
$ws.Range("G4").Value = "7"
$ws.Range("G5").Value = "7"
$ws.Range("G6").Value = "7"
$ws.Range("G7").Value = "7"
$ws.Range("G8").Value = "7"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").Value = "7"
$ws.Range("G10").Value = "7"
$ws.Range("G11").Value = "7"
$ws.Range("G12").Value = "7"
$ws.Range("G13").Value = "7"
$ws.Range("G14").Value = "7"
$ws.Range("G15").Value = "7"
$ws.Range("G16").Value = "7"
$ws.Range("G17").Value = "7"
$ws.Range("G18").Value = "7"
$ws.Range("G19").Value = "7"
$ws.Range("G20").Value = "7"
$ws.Range("G21").Value = "7"
$ws.Range("G22").Value = "7"
$ws.Range("G23").Value = "7"
$ws.Range("G24").Value = "7"
$ws.Range("G25").Value = "7"
$ws.Range("G26").Value = "7"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("G27").Value = "7"
$ws.Range("G28").Value = "7"
$ws.Range("G29").Value = "7"
$ws.Range("G30").Value = "7"
$ws.Range("G31").Value = "7"
$ws.Range("G32").Value = "7"
$ws.Range("G33").Value = "7"
$ws.Range("G34").Value = "7"
$ws.Range("G35").Value = "7"
$ws.Range("G36").Value = "7"
$ws.Range("G37").Value = "7"
$ws.Range("G38").Value = "7"
$ws.Range("G39").Value = "7"
$ws.Range("G40").Value = "7"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "7"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "7"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "7"
$ws.Range("G44").Value = "7"
$ws.Range("G45").Value = "7"
$ws.Range("G46").Value = "7"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "7"
$ws.Range("G48").Value = "7"
$ws.Range("G49").Value = "7"
$ws.Range("G50").Value = "7"
$ws.Range("G51").Value = "7"
